$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 274.80356
$ws.Cells.Item(17, 10).Value = 277.9818
$ws.Cells.Item(17, 12).Value = 833.9454000000001
$ws.Cells.Item(17, 14).Value = -1169.9454
$ws.Cells.Item(86, 8).Value = 2721.5
$ws.Cells.Item(86, 9).Value = 2490.3
$ws.Cells.Item(86, 10).Value = 3010.5
$ws.Cells.Item(86, 11).Value = 2490.3
$ws.Cells.Item(86, 12).Value = 3010.5
$ws.Cells.Item(86, 13).Value = -1367.3
$ws.Cells.Item(86, 14).Value = -5256.5
$ws.Cells.Item(89, 8).Value = 2721.5
$ws.Cells.Item(89, 9).Value = 2490.3
$ws.Cells.Item(89, 10).Value = 3010.5
$ws.Cells.Item(89, 11).Value = 12451.5
$ws.Cells.Item(89, 12).Value = 15052.5
$ws.Cells.Item(89, 13).Value = -6835.5
$ws.Cells.Item(89, 14).Value = -26284.5
$ws.Cells.Item(92, 8).Value = 538.5
$ws.Cells.Item(92, 9).Value = 551.3333
$ws.Cells.Item(92, 10).Value = 500
$ws.Cells.Item(92, 11).Value = 551.3333
$ws.Cells.Item(92, 12).Value = 500
$ws.Cells.Item(92, 13).Value = 696.6667
$ws.Cells.Item(92, 14).Value = -2996
$ws.Cells.Item(100, 8).Value = 1449.95
$ws.Cells.Item(100, 9).Value = 1166.5834
$ws.Cells.Item(100, 10).Value = 1875
$ws.Cells.Item(100, 11).Value = 1166.5834
$ws.Cells.Item(100, 12).Value = 1875
$ws.Cells.Item(100, 13).Value = -625.5834
$ws.Cells.Item(100, 14).Value = -2957

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 80001
$ws.Cells.Item(45, 9).Value = 202979.8
$ws.Cells.Item(45, 10).Value = 3139.25
$ws.Cells.Item(45, 11).Value = 202979.8
$ws.Cells.Item(45, 12).Value = 3139.25
$ws.Cells.Item(45, 13).Value = -202602.8
$ws.Cells.Item(45, 14).Value = -3893.25
$ws.Cells.Item(97, 8).Value = 399.66666
$ws.Cells.Item(97, 9).Value = 352.33334
$ws.Cells.Item(97, 10).Value = 525.8889
$ws.Cells.Item(97, 11).Value = 352.33334
$ws.Cells.Item(97, 12).Value = 525.8889
$ws.Cells.Item(97, 13).Value = 143.66666
$ws.Cells.Item(97, 14).Value = -1517.8889
$ws.Cells.Item(102, 8).Value = 1885.3846
$ws.Cells.Item(102, 9).Value = 1773.6364
$ws.Cells.Item(102, 10).Value = 2500
$ws.Cells.Item(102, 11).Value = 1773.6364
$ws.Cells.Item(102, 12).Value = 2500
$ws.Cells.Item(102, 13).Value = -151.6364000000001
$ws.Cells.Item(102, 14).Value = -5744
$ws.Cells.Item(122, 8).Value = 1853
$ws.Cells.Item(122, 9).Value = 1299.5
$ws.Cells.Item(122, 10).Value = 2960
$ws.Cells.Item(122, 11).Value = 3898.5
$ws.Cells.Item(122, 12).Value = 8880
$ws.Cells.Item(122, 13).Value = -1448.5
$ws.Cells.Item(122, 14).Value = -13780

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1909.091
$ws.Cells.Item(86, 9).Value = 1588.2354
$ws.Cells.Item(86, 11).Value = 1588.2354
$ws.Cells.Item(86, 13).Value = -465.2354
$ws.Cells.Item(89, 8).Value = 1909.091
$ws.Cells.Item(89, 9).Value = 1588.2354
$ws.Cells.Item(89, 11).Value = 7941.177
$ws.Cells.Item(89, 13).Value = -2325.177
$ws.Cells.Item(94, 8).Value = 989.4167
$ws.Cells.Item(94, 9).Value = 882.7143
$ws.Cells.Item(94, 10).Value = 1736.3334
$ws.Cells.Item(94, 11).Value = 882.7143
$ws.Cells.Item(94, 12).Value = 1736.3334
$ws.Cells.Item(94, 13).Value = -431.7143
$ws.Cells.Item(94, 14).Value = -2638.3334
$ws.Cells.Item(99, 8).Value = 833.8
$ws.Cells.Item(99, 9).Value = 757.8333
$ws.Cells.Item(99, 10).Value = 947.75
$ws.Cells.Item(99, 11).Value = 757.8333
$ws.Cells.Item(99, 12).Value = 947.75
$ws.Cells.Item(99, 13).Value = 740.1667
$ws.Cells.Item(99, 14).Value = -3943.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 1269.75
$ws.Cells.Item(7, 9).Value = 2012.6
$ws.Cells.Item(7, 10).Value = 31.666666
$ws.Cells.Item(7, 11).Value = 2012.6
$ws.Cells.Item(7, 12).Value = 31.666666
$ws.Cells.Item(7, 13).Value = -1899.6
$ws.Cells.Item(7, 14).Value = -257.666666
$ws.Cells.Item(105, 8).Value = 1074.75
$ws.Cells.Item(105, 9).Value = 911
$ws.Cells.Item(105, 10).Value = 1566
$ws.Cells.Item(105, 11).Value = 911
$ws.Cells.Item(105, 12).Value = 1566
$ws.Cells.Item(105, 13).Value = 836
$ws.Cells.Item(105, 14).Value = -5060
$ws.Cells.Item(132, 8).Value = 6138.25
$ws.Cells.Item(132, 9).Value = 8382.4
$ws.Cells.Item(132, 10).Value = 4535.2856
$ws.Cells.Item(132, 11).Value = 25147.2
$ws.Cells.Item(132, 12).Value = 13605.8568
$ws.Cells.Item(132, 13).Value = -22617.2
$ws.Cells.Item(132, 14).Value = -18665.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1123.96
$ws.Cells.Item(68, 9).Value = 772.5714
$ws.Cells.Item(68, 10).Value = 1571.1818
$ws.Cells.Item(68, 11).Value = 2317.7142
$ws.Cells.Item(68, 12).Value = 4713.5454
$ws.Cells.Item(68, 13).Value = -1506.7142
$ws.Cells.Item(68, 14).Value = -6335.5454
$ws.Cells.Item(71, 8).Value = 1123.96
$ws.Cells.Item(71, 9).Value = 772.5714
$ws.Cells.Item(71, 10).Value = 1571.1818
$ws.Cells.Item(71, 11).Value = 6953.1426
$ws.Cells.Item(71, 12).Value = 14140.6362
$ws.Cells.Item(71, 13).Value = -2897.1426
$ws.Cells.Item(71, 14).Value = -22252.6362
$ws.Cells.Item(113, 8).Value = 445.94287
$ws.Cells.Item(113, 10).Value = 424.82144
$ws.Cells.Item(113, 12).Value = 1274.46432
$ws.Cells.Item(113, 14).Value = -5614.46432
$ws.Cells.Item(129, 8).Value = 816.58826
$ws.Cells.Item(129, 9).Value = 272.375
$ws.Cells.Item(129, 10).Value = 1300.3334
$ws.Cells.Item(129, 11).Value = 817.125
$ws.Cells.Item(129, 12).Value = 3901.0002
$ws.Cells.Item(129, 13).Value = 4182.875
$ws.Cells.Item(129, 14).Value = -13901.0002
$ws.Cells.Item(131, 8).Value = 1164466.4
$ws.Cells.Item(131, 10).Value = 1317504.9
$ws.Cells.Item(131, 12).Value = 3952514.7
$ws.Cells.Item(131, 14).Value = -3962594.7
$ws.Cells.Item(136, 8).Value = 2073.5386
$ws.Cells.Item(136, 9).Value = 1494.75
$ws.Cells.Item(136, 10).Value = 2999.6
$ws.Cells.Item(136, 11).Value = 4484.25
$ws.Cells.Item(136, 12).Value = 8998.799999999999
$ws.Cells.Item(136, 13).Value = 615.75
$ws.Cells.Item(136, 14).Value = -19198.8
$ws.Cells.Item(138, 8).Value = 3398
$ws.Cells.Item(138, 9).Value = 3568.75
$ws.Cells.Item(138, 10).Value = 2942.6667
$ws.Cells.Item(138, 11).Value = 10706.25
$ws.Cells.Item(138, 12).Value = 8828.000100000001
$ws.Cells.Item(138, 13).Value = -5566.25
$ws.Cells.Item(138, 14).Value = -19108.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1178.75
$ws.Cells.Item(97, 9).Value = 1220.8334
$ws.Cells.Item(97, 10).Value = 1052.5
$ws.Cells.Item(97, 11).Value = 1220.8334
$ws.Cells.Item(97, 12).Value = 1052.5
$ws.Cells.Item(97, 13).Value = -724.8334
$ws.Cells.Item(97, 14).Value = -2044.5
$ws.Cells.Item(102, 8).Value = 1286.3636
$ws.Cells.Item(102, 9).Value = 1319.6364
$ws.Cells.Item(102, 10).Value = 1253.091
$ws.Cells.Item(102, 11).Value = 1319.6364
$ws.Cells.Item(102, 12).Value = 1253.091
$ws.Cells.Item(102, 13).Value = 302.3635999999999
$ws.Cells.Item(102, 14).Value = -4497.091
$ws.Cells.Item(122, 8).Value = 4564.6665
$ws.Cells.Item(122, 9).Value = 4478.7144
$ws.Cells.Item(122, 10).Value = 5166.3335
$ws.Cells.Item(122, 11).Value = 13436.1432
$ws.Cells.Item(122, 12).Value = 15499.0005
$ws.Cells.Item(122, 13).Value = -10986.1432
$ws.Cells.Item(122, 14).Value = -20399.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(30, 8).Value = 5010.6665
$ws.Cells.Item(30, 9).Value = 5010.6665
$ws.Cells.Item(30, 11).Value = 5010.6665
$ws.Cells.Item(30, 13).Value = -4902.6665
$ws.Cells.Item(32, 8).Value = 999
$ws.Cells.Item(32, 9).Value = 999
$ws.Cells.Item(32, 11).Value = 999
$ws.Cells.Item(32, 13).Value = -682
$ws.Cells.Item(40, 8).Value = 3625.05
$ws.Cells.Item(40, 9).Value = 3681.3125
$ws.Cells.Item(40, 10).Value = 3400
$ws.Cells.Item(40, 11).Value = 3681.3125
$ws.Cells.Item(40, 12).Value = 3400
$ws.Cells.Item(40, 13).Value = -3545.3125
$ws.Cells.Item(40, 14).Value = -3672
$ws.Cells.Item(93, 8).Value = 1989.7667
$ws.Cells.Item(93, 9).Value = 1727.2858
$ws.Cells.Item(93, 10).Value = 2219.4375
$ws.Cells.Item(93, 11).Value = 1727.2858
$ws.Cells.Item(93, 12).Value = 2219.4375
$ws.Cells.Item(93, 13).Value = -479.2858000000001
$ws.Cells.Item(93, 14).Value = -4715.4375
$ws.Cells.Item(122, 8).Value = 8405.632
$ws.Cells.Item(122, 9).Value = 15148.25
$ws.Cells.Item(122, 11).Value = 45444.75
$ws.Cells.Item(122, 13).Value = -42994.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 10011911
$ws.Cells.Item(96, 10).Value = 16658.572
$ws.Cells.Item(96, 12).Value = 16658.572
$ws.Cells.Item(96, 14).Value = -19404.572
$ws.Cells.Item(122, 8).Value = 1689.3889
$ws.Cells.Item(122, 9).Value = 1446.2727
$ws.Cells.Item(122, 10).Value = 2071.4285
$ws.Cells.Item(122, 11).Value = 4338.8181
$ws.Cells.Item(122, 12).Value = 6214.2855
$ws.Cells.Item(122, 13).Value = -1888.8181
$ws.Cells.Item(122, 14).Value = -11114.2855
